# Update countries & provincias Spain
# Refresh COVID-19 stats for several countries and fix the row order for
# "Singapur" (now ahead of "Bielorrusia") and "Argelia" (now ahead of "Marruecos").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 2).Value = 653397
$ws.Cells.Item(4, 3).Value = 5249
$ws.Cells.Item(4, 4).Value = 56618
$ws.Cells.Item(4, 5).Value = 563374
$ws.Cells.Item(4, 6).Value = 13516
$ws.Cells.Item(4, 7).Value = 817
$ws.Cells.Item(4, 8).Value = 33405

# Row 15
$ws.Cells.Item(15, 6).Value = 6634

# Row 18
$ws.Cells.Item(18, 5).Value = 10051
$ws.Cells.Item(18, 7).Value = 42
$ws.Cells.Item(18, 8).Value = 1281

# Row 34
$ws.Cells.Item(34, 2).Value = 6848
$ws.Cells.Item(34, 3).Value = 51
$ws.Cells.Item(34, 5).Value = 6664
$ws.Cells.Item(34, 7).Value = 2
$ws.Cells.Item(34, 8).Value = 152

# Row 44 -> Singapur
$ws.Cells.Item(44, 1).Value = "Singapur"
$ws.Cells.Item(44, 2).Value = 4427
$ws.Cells.Item(44, 3).Value = 728
$ws.Cells.Item(44, 4).Value = 683
$ws.Cells.Item(44, 5).Value = 3734
$ws.Cells.Item(44, 6).Value = 29
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 10

# Row 45 -> Bielorrusia
$ws.Cells.Item(45, 1).Value = "Bielorrusia"
$ws.Cells.Item(45, 2).Value = 4204
$ws.Cells.Item(45, 3).Value = 476
$ws.Cells.Item(45, 4).Value = 203
$ws.Cells.Item(45, 5).Value = 3961
$ws.Cells.Item(45, 6).Value = 65
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 8).Value = 40

# Row 46 -> Ucrania
$ws.Cells.Item(46, 1).Value = "Ucrania"
$ws.Cells.Item(46, 2).Value = 4161
$ws.Cells.Item(46, 3).Value = 397
$ws.Cells.Item(46, 4).Value = 186
$ws.Cells.Item(46, 5).Value = 3859
$ws.Cells.Item(46, 6).Value = 45
$ws.Cells.Item(46, 7).Value = 8
$ws.Cells.Item(46, 8).Value = 116

# Row 47 -> Catar
$ws.Cells.Item(47, 1).Value = "Catar"
$ws.Cells.Item(47, 2).Value = 4103
$ws.Cells.Item(47, 3).Value = 392
$ws.Cells.Item(47, 4).Value = 415
$ws.Cells.Item(47, 5).Value = 3681
$ws.Cells.Item(47, 6).Value = 37
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 7

# Row 48 -> Republica Dominicana
$ws.Cells.Item(48, 1).Value = "Republica Dominicana"
$ws.Cells.Item(48, 2).Value = 3755
$ws.Cells.Item(48, 3).Value = 141
$ws.Cells.Item(48, 4).Value = 215
$ws.Cells.Item(48, 5).Value = 3344
$ws.Cells.Item(48, 6).Value = 143
$ws.Cells.Item(48, 7).Value = 7
$ws.Cells.Item(48, 8).Value = 196

# Row 49 -> Panama
$ws.Cells.Item(49, 1).Value = "Panama"
$ws.Cells.Item(49, 2).Value = 3751
$ws.Cells.Item(49, 4).Value = 75
$ws.Cells.Item(49, 5).Value = 3573
$ws.Cells.Item(49, 6).Value = 106
$ws.Cells.Item(49, 8).Value = 103

# Row 57 -> Argelia
$ws.Cells.Item(57, 1).Value = "Argelia"
$ws.Cells.Item(57, 2).Value = 2268
$ws.Cells.Item(57, 3).Value = 108
$ws.Cells.Item(57, 4).Value = 783
$ws.Cells.Item(57, 5).Value = 1137
$ws.Cells.Item(57, 6).Value = 60
$ws.Cells.Item(57, 7).Value = 12
$ws.Cells.Item(57, 8).Value = 348

# Row 58 -> Marruecos
$ws.Cells.Item(58, 1).Value = "Marruecos"
$ws.Cells.Item(58, 2).Value = 2251
$ws.Cells.Item(58, 3).Value = 227
$ws.Cells.Item(58, 4).Value = 247
$ws.Cells.Item(58, 5).Value = 1876
$ws.Cells.Item(58, 6).Value = 1
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 128

# Row 59 -> Grecia
$ws.Cells.Item(59, 1).Value = "Grecia"
$ws.Cells.Item(59, 2).Value = 2207
$ws.Cells.Item(59, 3).Value = 15
$ws.Cells.Item(59, 4).Value = 269
$ws.Cells.Item(59, 5).Value = 1833
$ws.Cells.Item(59, 6).Value = 69
$ws.Cells.Item(59, 7).Value = 3
$ws.Cells.Item(59, 8).Value = 105

# Row 68
$ws.Cells.Item(68, 2).Value = 1434
$ws.Cells.Item(68, 3).Value = 19
$ws.Cells.Item(68, 4).Value = 856
$ws.Cells.Item(68, 5).Value = 498
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = 80

# Row 110
$ws.Cells.Item(110, 2).Value = 340
$ws.Cells.Item(110, 3).Value = 34
$ws.Cells.Item(110, 4).Value = 76
$ws.Cells.Item(110, 5).Value = 261
